$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @(
    'Blouse,Kaftan'
    'Kaftan,Tee'
    'Kaftan,Tee'
    'Jumpsuit,Tee'
    'Kaftan,Jumpsuit'
    'Halter,Tee'
    'Jumpsuit,Blouse'
    'Halter,Jumpsuit'
    'Tee,Halter'
    'Jumpsuit,Tee'
    'Jumpsuit,Dress'
    'Trunks,Jumpsuit'
    'Jodhpurs,Tee'
    'Jumpsuit,Halter'
    'Jumpsuit,Kaftan'
    'Jumpsuit,Blouse'
    'Parka,Jumpsuit'
    'Jumpsuit,Tee'
    'Jumpsuit,Dress'
    'Blazer,Jumpsuit'
    'Blouse,Jumpsuit'
    'Jumpsuit,Kaftan'
    'Halter,Jumpsuit'
    'Blazer,Halter'
    'Jumpsuit,Blazer'
    'Jumpsuit,Halter'
    'Halter,Blazer'
    'Jumpsuit,Halter'
    'Parka,Tee'
    'Trunks,Caftan'
    'Caftan,Jumpsuit'
    'Halter,Blazer'
    'Blazer,Jumpsuit'
    'Blouse,Blazer'
    'Jumpsuit,Trunks'
    'Blazer,Halter'
    'Trunks,Cutoffs'
    'Caftan,Trunks'
    'Blazer,Blouse'
    'Jumpsuit,Blazer'
    'Jumpsuit,Tee'
    'Halter,Blazer'
    'Jodhpurs,Kaftan'
    'Jumpsuit,Halter'
    'Halter,Blazer'
    'Halter,Blazer'
    'Jumpsuit,Tee'
    'Halter,Caftan'
    'Jumpsuit,Dress'
    'Jumpsuit,Halter'
    'Trunks,Jumpsuit'
    'Jumpsuit,Blouse'
    'Halter,Blazer'
    'Jumpsuit,Blazer'
    'Jumpsuit,Dress'
    'Jumpsuit,Dress'
    'Jumpsuit,Blouse'
    'Trunks,Jumpsuit'
    'Blouse,Jumpsuit'
    'Jumpsuit,Dress'
    'Tee,Kaftan'
    'Jumpsuit,Kaftan'
    'Tee,Kaftan'
    'Kaftan,Tee'
    'Jumpsuit,Tee'
    'Kaftan,Jodhpurs'
    'Dress,Jumpsuit'
    'Kaftan,Tee'
    'Parka,Jumpsuit'
    'Parka,Jumpsuit'
    'Jumpsuit,Parka'
    'Jumpsuit,Dress'
    'Jumpsuit,Tee'
    'Caftan,Parka'
    'Jumpsuit,Dress'
    'Jumpsuit,Kaftan'
    'Jumpsuit,Kaftan'
    'Jumpsuit,Trunks'
    'Jumpsuit,Halter'
    'Jumpsuit,Jeggings'
    'Jumpsuit,Caftan'
    'Jumpsuit,Halter'
    'Halter,Blazer'
    'Blazer,Jeggings'
    'Blazer,Halter'
    'Tee,Kaftan'
    'Jumpsuit,Blouse'
    'Jumpsuit,Trunks'
    'Halter,Tee'
    'Parka,Sweatpants'
    'Jumpsuit,Blouse'
    'Halter,Jumpsuit'
    'Halter,Blouse'
    'Blazer,Halter'
    'Blouse,Jumpsuit'
    'Jumpsuit,Caftan'
    'Parka,Blouse'
    'Jumpsuit,Dress'
    'Sweatpants,Jumpsuit'
    'Parka,Jumpsuit'
    'Jumpsuit,Dress'
    'Parka,Blouse'
    'Halter,Kaftan'
    'Tee,Kaftan'
    'Jumpsuit,Blouse'
    'Halter,Jumpsuit'
    'Jumpsuit,Blouse'
    'Jumpsuit,Kaftan'
    'Jodhpurs,Trunks'
    'Gauchos,Sweatpants'
    'Parka,Caftan'
    'Jumpsuit,Blouse'
    'Jumpsuit,Blouse'
    'Halter,Jodhpurs'
    'Tee,Kaftan'
    'Trunks,Jodhpurs'
    'Caftan,Jumpsuit'
    'Halter,Tee'
    'Caftan,Jumpsuit'
    'Trunks,Jumpsuit'
    'Trunks,Caftan'
    'Jumpsuit,Blouse'
    'Trunks,Parka'
    'Jumpsuit,Tee'
    'Caftan,Jumpsuit'
    'Jumpsuit,Tee'
    'Caftan,Blazer'
    'Jumpsuit,Blouse'
    'Halter,Jumpsuit'
    'Jumpsuit,Blouse'
    'Jumpsuit,Kaftan'
    'Jumpsuit,Dress'
    'Tee,Kaftan'
    'Jumpsuit,Kaftan'
    'Halter,Jumpsuit'
    'Jumpsuit,Dress'
    'Jumpsuit,Trunks'
    'Jumpsuit,Halter'
    'Jumpsuit,Cutoffs'
    'Jumpsuit,Caftan'
    'Kaftan,Tee'
    'Halter,Jumpsuit'
    'Cutoffs,Jodhpurs'
    'Halter,Blazer'
    'Jumpsuit,Kaftan'
    'Tee,Kaftan'
    'Jumpsuit,Halter'
    'Halter,Tee'
    'Cutoffs,Trunks'
    'Caftan,Blouse'
    'Jumpsuit,Hoodie'
    'Cutoffs,Jumpsuit'
    'Parka,Gauchos'
    'Trunks,Blazer'
    'Jumpsuit,Blouse'
    'Cutoffs,Jodhpurs'
    'Jumpsuit,Dress'
    'Cutoffs,Jumpsuit'
    'Tee,Kaftan'
    'Kaftan,Tee'
    'Jodhpurs,Kaftan'
    'Blazer,Halter'
    'Tee,Kaftan'
    'Halter,Tee'
    'Jumpsuit,Blouse'
    'Jumpsuit,Halter'
    'Jumpsuit,Kaftan'
    'Trunks,Caftan'
    'Jumpsuit,Kaftan'
    'Blazer,Halter'
    'Blouse,Parka'
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $values[$i]
}
